$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 58:59 (old rows 58-61 shift down to 60-63)
$ws.Rows("58:59").Insert()

# New row 58: R26 pulldown resistor
$ws.Range("A58").Value = "R26"
$ws.Range("B58").Value = "0R"
$ws.Range("C58").Value = "RESISTOR0603"
$ws.Range("D58").Value = "0603-RES"
$ws.Range("E58").Value = "Resistor"

# New row 59: R27 pulldown resistor
$ws.Range("A59").Value = "R27"
$ws.Range("B59").Value = "50K"
$ws.Range("C59").Value = "RESISTOR0603"
$ws.Range("D59").Value = "0603-RES"
$ws.Range("E59").Value = "Resistor"

# Existing R30 row (now at row 60) changes Value from 0R to 50K
$ws.Range("B60").Value = "50K"

# Update the defined name range to cover the two extra rows
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!ww101board") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$63"
    }
}

# Update selection to match final state
$ws.Range("B61").Select() | Out-Null
